$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1
$ws.Cells.Item($newRow, 1).Value = "IM2325000001"
